{"js": "// Apply the diff: update the date line and 25 multiplication-table\n// answers throughout the document body. Every \"old\" string occurs\n// exactly once in this document, so a direct search + replace on\n// each pair is unambiguous and idempotent.\nconst replacements = [\n  ['2024-09-17 Tuesday', '2024-09-18 Wednesday'],\n  ['761\u00d75=3805', '114\u00d73=342'],\n  ['230\u00d73=690', '685\u00d76=4110'],\n  ['367\u00d75=1835', '106\u00d79=954'],\n  ['844\u00d79=7596', '530\u00d77=3710'],\n  ['418\u00d72=836', '744\u00d76=4464'],\n  ['181\u00d75=905', '661\u00d75=3305'],\n  ['261\u00d73=783', '680\u00d78=5440'],\n  ['423\u00d76=2538', '494\u00d75=2470'],\n  ['721\u00d73=2163', '127\u00d72=254'],\n  ['935\u00d72=1870', '315\u00d76=1890'],\n  ['698\u00d72=1396', '278\u00d73=834'],\n  ['484\u00d78=3872', '800\u00d76=4800'],\n  ['513\u00d72=1026', '736\u00d77=5152'],\n  ['814\u00d74=3256', '197\u00d72=394'],\n  ['620\u00d79=5580', '447\u00d77=3129'],\n  ['577\u00d76=3462', '751\u00d73=2253'],\n  ['944\u00d72=1888', '799\u00d75=3995'],\n  ['609\u00d75=3045', '675\u00d72=1350'],\n  ['144\u00d73=432', '782\u00d78=6256'],\n  ['199\u00d73=597', '455\u00d77=3185'],\n  ['612\u00d75=3060', '503\u00d78=4024'],\n  ['138\u00d76=828', '436\u00d73=1308'],\n  ['891\u00d77=6237', '365\u00d78=2920'],\n  ['759\u00d79=6831', '760\u00d78=6080'],\n  ['239\u00d76=1434', '899\u00d73=2697'],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true });\n  found.load('items');\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the diff: update the date line and 25 multiplication-table\n# answers throughout the document body. Every \"Find\" string occurs\n# exactly once in this document, so Find/Replace with MatchWholeWord\n# off but MatchCase on (and no wildcards) is unambiguous per pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"2024-09-17 Tuesday\"; Replace = \"2024-09-18 Wednesday\" },\n    @{ Find = \"761\u00d75=3805\"; Replace = \"114\u00d73=342\" },\n    @{ Find = \"230\u00d73=690\"; Replace = \"685\u00d76=4110\" },\n    @{ Find = \"367\u00d75=1835\"; Replace = \"106\u00d79=954\" },\n    @{ Find = \"844\u00d79=7596\"; Replace = \"530\u00d77=3710\" },\n    @{ Find = \"418\u00d72=836\"; Replace = \"744\u00d76=4464\" },\n    @{ Find = \"181\u00d75=905\"; Replace = \"661\u00d75=3305\" },\n    @{ Find = \"261\u00d73=783\"; Replace = \"680\u00d78=5440\" },\n    @{ Find = \"423\u00d76=2538\"; Replace = \"494\u00d75=2470\" },\n    @{ Find = \"721\u00d73=2163\"; Replace = \"127\u00d72=254\" },\n    @{ Find = \"935\u00d72=1870\"; Replace = \"315\u00d76=1890\" },\n    @{ Find = \"698\u00d72=1396\"; Replace = \"278\u00d73=834\" },\n    @{ Find = \"484\u00d78=3872\"; Replace = \"800\u00d76=4800\" },\n    @{ Find = \"513\u00d72=1026\"; Replace = \"736\u00d77=5152\" },\n    @{ Find = \"814\u00d74=3256\"; Replace = \"197\u00d72=394\" },\n    @{ Find = \"620\u00d79=5580\"; Replace = \"447\u00d77=3129\" },\n    @{ Find = \"577\u00d76=3462\"; Replace = \"751\u00d73=2253\" },\n    @{ Find = \"944\u00d72=1888\"; Replace = \"799\u00d75=3995\" },\n    @{ Find = \"609\u00d75=3045\"; Replace = \"675\u00d72=1350\" },\n    @{ Find = \"144\u00d73=432\"; Replace = \"782\u00d78=6256\" },\n    @{ Find = \"199\u00d73=597\"; Replace = \"455\u00d77=3185\" },\n    @{ Find = \"612\u00d75=3060\"; Replace = \"503\u00d78=4024\" },\n    @{ Find = \"138\u00d76=828\"; Replace = \"436\u00d73=1308\" },\n    @{ Find = \"891\u00d77=6237\"; Replace = \"365\u00d78=2920\" },\n    @{ Find = \"759\u00d79=6831\"; Replace = \"760\u00d78=6080\" },\n    @{ Find = \"239\u00d76=1434\"; Replace = \"899\u00d73=2697\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace)\n    # wdReplaceAll = 2, wdFindStop = 0\n    $find.Execute($r.Find, $true, $false, $false, $false, $false, $true, 0, $false, $r.Replace, 2) | Out-Null\n}\n"}
